# Appends new daily OHLCV rows (118-135) to the UNITDSPR.NS stock history
# sheet and backfills column R ("backup") on the previously-last rows
# (108-117) from blank to an explicit 0, matching the "break out stock.yaml
# completed" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rows 108-117: column R ("backup") goes from blank to numeric 0 ---
for ($r = 108; $r -le 117; $r++) {
    $ws.Cells.Item($r, 18).Value = 0   # R = column 18
}

# --- 2) New rows 118-135: full OHLCV + date-part rows ---
# Columns: row, A Datetime(serial), B Open, C High, D Low, E Close,
#          F AdjClose, G Volume, H Year, I Month, J Day, K Hour,
#          L Minute, M Second, N Week, hasOPQ (O/P/Q filled with 0)
$newRows = @(
    @(118,45617,1500,1504.949951171875,1481.550048828125,1493.050048828125,1493.050048828125,545772,2024,11,21,0,0,0,47,0),
    @(119,45618,1493,1508,1489.699951171875,1500.150024414062,1500.150024414062,621662,2024,11,22,0,0,0,47,0),
    @(120,45621,1518.900024414062,1543.699951171875,1476.150024414062,1486.099975585938,1486.099975585938,3600031,2024,11,25,0,0,0,48,0),
    @(121,45622,1486.099975585938,1528.75,1466.900024414062,1511.150024414062,1511.150024414062,658663,2024,11,26,0,0,0,48,0),
    @(122,45623,1518.699951171875,1526.949951171875,1500.050048828125,1515.75,1515.75,547190,2024,11,27,0,0,0,48,0),
    @(123,45624,1521.400024414062,1522.75,1495.449951171875,1502.75,1502.75,566858,2024,11,28,0,0,0,48,0),
    @(124,45625,1516,1541.949951171875,1501.150024414062,1529.099975585938,1529.099975585938,1224425,2024,11,29,0,0,0,48,0),
    @(125,45628,1542,1552.849975585938,1520,1532.199951171875,1532.199951171875,554366,2024,12,2,0,0,0,49,0),
    @(126,45629,1532.400024414062,1545,1513.099975585938,1542.949951171875,1542.949951171875,522548,2024,12,3,0,0,0,49,1),
    @(127,45630,1549,1549,1520.650024414062,1526.050048828125,1526.050048828125,618463,2024,12,4,0,0,0,49,1),
    @(128,45631,1526.199951171875,1535,1515.150024414062,1530.099975585938,1530.099975585938,423566,2024,12,5,0,0,0,49,1),
    @(129,45632,1537.75,1537.75,1513.050048828125,1516.599975585938,1516.599975585938,323615,2024,12,6,0,0,0,49,1),
    @(130,45635,1512.699951171875,1522.650024414062,1492.25,1506.25,1506.25,516426,2024,12,9,0,0,0,50,1),
    @(131,45636,1517,1526.150024414062,1500.050048828125,1513.5,1513.5,522996,2024,12,10,0,0,0,50,1),
    @(132,45637,1520,1522.199951171875,1506.75,1516.800048828125,1516.800048828125,324261,2024,12,11,0,0,0,50,1),
    @(133,45638,1510.800048828125,1532.75,1510.800048828125,1525.949951171875,1525.949951171875,540314,2024,12,12,0,0,0,50,1),
    @(134,45639,1528.150024414062,1529,1482.849975585938,1512.050048828125,1512.050048828125,462501,2024,12,13,0,0,0,50,1),
    @(135,45642,1512.050048828125,1544.300048828125,1467.349975585938,1519.199951171875,1519.199951171875,805772,2024,12,16,0,0,0,51,1)
)

# Column A keeps the same datetime number format used by the existing rows
$ws.Range("A118:A135").NumberFormat = "YYYY-MM-DD HH:MM:SS"

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]    # A Datetime
    $ws.Cells.Item($r, 2).Value = $row[2]    # B Open
    $ws.Cells.Item($r, 3).Value = $row[3]    # C High
    $ws.Cells.Item($r, 4).Value = $row[4]    # D Low
    $ws.Cells.Item($r, 5).Value = $row[5]    # E Close
    $ws.Cells.Item($r, 6).Value = $row[6]    # F Adj Close
    $ws.Cells.Item($r, 7).Value = $row[7]    # G Volume
    $ws.Cells.Item($r, 8).Value = $row[8]    # H Year
    $ws.Cells.Item($r, 9).Value = $row[9]    # I Month
    $ws.Cells.Item($r, 10).Value = $row[10]  # J Day
    $ws.Cells.Item($r, 11).Value = $row[11]  # K Hour
    $ws.Cells.Item($r, 12).Value = $row[12]  # L Minute
    $ws.Cells.Item($r, 13).Value = $row[13]  # M Second
    $ws.Cells.Item($r, 14).Value = $row[14]  # N Week

    if ($row[15] -eq 1) {
        $ws.Cells.Item($r, 15).Value = 0     # O isPivot
        $ws.Cells.Item($r, 16).Value = 0     # P two_line_structure
        $ws.Cells.Item($r, 17).Value = 0     # Q detect_structure
    }
    # R (backup) is left blank for every new row, matching the source diff.
}
